# Update cryptocurrency price (column D) and hourly volume change (column E)
# values to reflect the latest scrape, per the Sun Sep 29 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.690.32'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.649.75'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.86'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.95'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.639'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.69%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.397'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.81'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.41%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.54'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.39%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.123.04'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.529.97'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.662.56'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.52'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.75'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '349.14'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.42'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.20%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.44'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.79'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +8.74%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.55'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.60'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '561.44'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +6.71%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.02'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.86%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.78'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.51'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.49'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.35'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '155.25'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '160.45'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.06'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0603'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.49%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.54'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.636'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0255'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.85%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.64'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0244'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.32%  '
